$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Naming Rule")
$ws2 = $wb.Worksheets.Item("Configuration Options")

# --- Naming Rule sheet ---
# Remove the existing hyperlink first so the row-insert below doesn't leave
# a stale/duplicated hyperlink entry pointing at the old address.
$ws1.Hyperlinks.Delete()

# Insert a new row at 11 (pushes the old row 12 - the "참고" / hyperlink row -
# down to row 13), matching the diff's new row B2:S13 dimension.
$ws1.Rows.Item(11).Insert()

# Re-create the blank "D12 / E12" placeholder cells (same spot as before the
# insert) using the formatting of the row that got pushed down to 13.
$ws1.Range("L13").Copy()
$ws1.Range("D12").PasteSpecial(-4122)
$ws1.Range("M13").Copy()
$ws1.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "Daebo" manufacturer option and "RSU GPS" / "Xsens GPS ROS (RTK)" /
# "Intel" options added by this commit.
$ws1.Range("H9").Value2 = "• D : Daebo"
$ws1.Range("N10").Value2 = "• RG : RSU GPS"
$ws1.Range("F11").Value2 = "• I : Intel"
$ws1.Range("N11").Value2 = "• XGR : Xsens GPS ROS (RTK)"

# Re-create the hyperlink on its new row (13) and restore its original cell
# formatting (Hyperlinks.Add resets the cell style as a side effect).
$ws1.Hyperlinks.Add($ws1.Range("M13"), "https://www.hanwhavision.com/ko/product-naming-rule/") | Out-Null
$ws1.Range("E13").Copy()
$ws1.Range("M13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Range("D20").Select() | Out-Null

# --- Configuration Options sheet ---
# New RSU GPS ROS model row.
$ws2.Range("C11").Value2 = "NO-CO2XGR"
$ws2.Range("D11").Value2 = "NVIDIA Orin RxAnt.2 Modem ROS"

# Ubuntu+Intel OBU model renamed from U-CO1XG / UBUNTU to UI-CO1XG / Ubuntu Intel OBU.
$ws2.Range("C5").Value2 = "UI-CO1XG"

# New Ubuntu+Intel RSU model row.
$ws2.Range("C12").Value2 = "UI-DR2RG"
$ws2.Range("D12").Value2 = "Ubuntu Intel RSU"

$ws2.Range("D5").Value2 = "Ubuntu Intel OBU"

$ws2.Range("D5").Select() | Out-Null
